$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0.9999999991123205"
$ws.Range("C2").Value = [double]"2.710008657200506e-10"
$ws.Range("D2").Value = [double]"1.339493385849563e-10"
$ws.Range("E2").Value = [double]"4.827292761601328e-10"
$ws.Range("F2").Value = 46073
